$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows at position 25 for new library.* keys, shifting
# blog.title..impact.sections.collaboration.description down to rows 35-91.
$ws.Rows("25:34").Insert()

# Fill newly inserted rows 25-34 with the new library.* translation keys.
$ws.Range("A25").Value = 'library.table_of_contents'
$ws.Range("B25").Value = 'Table of Contents'
$ws.Range("C25").Value = 'Mục Lục'

$ws.Range("A26").Value = 'library.key_takeaways'
$ws.Range("B26").Value = 'Key Takeaways'
$ws.Range("C26").Value = 'Điểm Chính'

$ws.Range("A27").Value = 'library.book_not_found'
$ws.Range("B27").Value = 'Book not found'
$ws.Range("C27").Value = 'Không tìm thấy sách'

$ws.Range("A28").Value = 'library.book_not_found_description'
$ws.Range("B28").Value = 'Sorry, we couldn''t find the book you''re looking for.'
$ws.Range("C28").Value = 'Xin lỗi, chúng tôi không thể tìm thấy cuốn sách bạn đang tìm kiếm.'

$ws.Range("A29").Value = 'library.reading_time'
$ws.Range("B29").Value = '{{minutes}} min read'
$ws.Range("C29").Value = '{{minutes}} phút đọc'

$ws.Range("A30").Value = 'library.chapter'
$ws.Range("B30").Value = 'Chapter {{number}}'
$ws.Range("C30").Value = 'Chương {{number}}'

$ws.Range("A31").Value = 'library.previous_chapter'
$ws.Range("B31").Value = 'Previous Chapter'
$ws.Range("C31").Value = 'Chương Trước'

$ws.Range("A32").Value = 'library.next_chapter'
$ws.Range("B32").Value = 'Next Chapter'
$ws.Range("C32").Value = 'Chương Tiếp Theo'

$ws.Range("A33").Value = 'library.by'
$ws.Range("B33").Value = 'by'
$ws.Range("C33").Value = 'bởi'

$ws.Range("A34").Value = 'library.last_updated'
$ws.Range("B34").Value = 'Last updated {{date}}'
$ws.Range("C34").Value = 'Cập nhật lần cuối {{date}}'

# Append two new common.* rows at the end (rows 92-93).
$ws.Range("A92").Value = 'common.back_to_library'
$ws.Range("B92").Value = 'Back to Library'
$ws.Range("C92").Value = 'Quay Lại Thư Viện'

$ws.Range("A93").Value = 'common.published_on'
$ws.Range("B93").Value = 'Published on {{date}}'
$ws.Range("C93").Value = 'Xuất bản ngày {{date}}'
